$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. First paragraph: "This is a Microsoft word document."
#    -> append two trailing spaces, then a red parenthetical note,
#    typed as three separate insertions (mirrors three distinct
#    runs sharing the same red-color formatting).
# ------------------------------------------------------------------
$para1 = $d.Paragraphs(1).Range

# Position right before the paragraph mark at the end of paragraph 1.
$insertPoint = $d.Range($para1.End - 1, $para1.End - 1)
$insertPoint.InsertAfter("  ")

$run1 = $d.Range($insertPoint.End, $insertPoint.End)
$run1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$run1.Font.Color = 255

$run2 = $d.Range($run1.End, $run1.End)
$run2.InsertAfter("rsion for main branch")
$run2.Font.Color = 255

$run3 = $d.Range($run2.End, $run2.End)
$run3.InsertAfter(")")
$run3.Font.Color = 255

# ------------------------------------------------------------------
# 2. Remove the trailing "ank God almighty, we are free at last."
#    paragraph entirely (it is the final paragraph in the document).
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastPara.Range.Delete()

# ------------------------------------------------------------------
# 3. Drop the handful of unused styles that came along for the ride
#    (none of them are referenced anywhere in the body text). Walk
#    the list back-to-front by collection position so each Delete()
#    never invalidates an index we still need to look up.
# ------------------------------------------------------------------
$unusedStyleNames = @(
    "podcast-toolssubscribe-links",
    "generic-title",
    "subscribe-more-info",
    "subscribe",
    "audio-tool",
    "Heading4Char",
    "Heading2Char",
    "Hyperlink",
    "apple-converted-space",
    "Heading4",
    "Heading2"
)
foreach ($styleName in $unusedStyleNames) {
    $style = $d.Styles($styleName)
    $style.Delete()
}
